$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template used to ship with 16 pre-filled sample rows (4-19). The
# exporter now adds rows to the sheet dynamically as needed, so the
# template only needs to keep its first two seed rows (4-5) with their
# formulas; the quantity seed values in S4/S5 are cleared too.
$ws.Range("S4:S5").ClearContents()

# Remove the now-unneeded extra template rows 6-19 (shifts rows below up).
$ws.Rows("6:19").Delete()

# Leave the cursor on K8, matching where the author left off editing.
$ws.Range("K8").Select()
